$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 7
$ws.Range("E3").Value = 9
$ws.Range("E4").Formula = "=SUM(E2:E3)"

$ws.Range("G7").Select() | Out-Null
